$d = $word.ActiveDocument

# 1. Replace "Enero" -> "Agosto" (first occurrence, the semester start month)
$d.Content.Find.Execute("Enero", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Agosto", 2)

# 2. Replace "Mayo" -> "Diciembre" (the semester end month)
$d.Content.Find.Execute("Mayo", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Diciembre", 2)
